$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Group" header in column I, matching the header style used
# by the rest of row 1 (copy format from H1).
$ws.Range("I1").Value = "Group"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# Fill in the Group value for the first data row; the second data row (I3)
# is intentionally left blank, matching the "required but not yet filled
# in" state the commit is capturing.
$ws.Range("I2").Value = "A"

# The existing data cells (A2:H3) pick up an explicit "apply font" style
# once the sheet is touched for the new required column.
$ws.Range("A2:H3").Style = "Normal"

# Leave the selection on I3, the next empty Group cell that now needs data.
$ws.Range("I3").Select()
